# Raul's Log - append new task-log entries (Jan 23-31, 2017 batch)
# and put sheet into "test mode" selection state per commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Row 1335
$ws.Cells.Item(1335, 1).Value = 'Pickup Small PA'
$ws.Cells.Item(1335, 2).Value = 42758
$ws.Cells.Item(1335, 3).Value = '1730'
$ws.Cells.Item(1335, 4).Value = 'SSB'
$ws.Cells.Item(1335, 5).Value = 'S126'
$ws.Cells.Item(1335, 6).Value = 'Return small PA with neck mic on cart to SSB N103'

# Row 1339
$ws.Cells.Item(1339, 1).Value = 'Pickup Small PA'
$ws.Cells.Item(1339, 2).Value = 42759
$ws.Cells.Item(1339, 3).Value = '1600'
$ws.Cells.Item(1339, 4).Value = 'DB'
$ws.Cells.Item(1339, 5).Value = '1014'
$ws.Cells.Item(1339, 6).Value = 'Return speaker, tripod, mixer, 3 mics with stands, one neck mic to DB 0003'
$ws.Rows.Item(1339).RowHeight = 30

# Row 1340
$ws.Cells.Item(1340, 1).Value = 'Pickup Projector'
$ws.Cells.Item(1340, 2).Value = 42759
$ws.Cells.Item(1340, 3).Value = '1730'
$ws.Cells.Item(1340, 4).Value = 'SSB'
$ws.Cells.Item(1340, 5).Value = 'W356'
$ws.Cells.Item(1340, 6).Value = 'Returjn doc camera on cart to SSB N103'

# Row 1341
$ws.Cells.Item(1341, 1).Value = 'Demo'
$ws.Cells.Item(1341, 2).Value = 42759
$ws.Cells.Item(1341, 3).Value = '1600'
$ws.Cells.Item(1341, 4).Value = 'DB'
$ws.Cells.Item(1341, 5).Value = '2008'

# Row 1342
$ws.Cells.Item(1342, 1).Value = 'Demo'
$ws.Cells.Item(1342, 2).Value = 42759
$ws.Cells.Item(1342, 3).Value = '1800'
$ws.Cells.Item(1342, 4).Value = 'OSG'
$ws.Cells.Item(1342, 5).Value = '2010'

# Row 1343
$ws.Cells.Item(1343, 1).Value = 'Demo'
$ws.Cells.Item(1343, 2).Value = 42759
$ws.Cells.Item(1343, 3).Value = '1900'
$ws.Cells.Item(1343, 4).Value = 'ACE'
$ws.Cells.Item(1343, 5).Value = '010'

# Row 1347
$ws.Cells.Item(1347, 1).Value = 'Setup PC'
$ws.Cells.Item(1347, 2).Value = 42760
$ws.Cells.Item(1347, 3).Value = '1630'
$ws.Cells.Item(1347, 4).Value = 'HNE'
$ws.Cells.Item(1347, 5).Value = '105'
$ws.Cells.Item(1347, 6).Value = 'Equipment from HNES 003'

# Row 1348
$ws.Cells.Item(1348, 1).Value = 'Setup Projector'
$ws.Cells.Item(1348, 2).Value = 42760
$ws.Cells.Item(1348, 3).Value = '1630'
$ws.Cells.Item(1348, 4).Value = 'HNE'
$ws.Cells.Item(1348, 5).Value = '105'
$ws.Cells.Item(1348, 6).Value = 'Equipment from HNES 003'

# Row 1349
$ws.Cells.Item(1349, 1).Value = 'Pickup Projector'
$ws.Cells.Item(1349, 2).Value = 42760
$ws.Cells.Item(1349, 3).Value = '1730'
$ws.Cells.Item(1349, 4).Value = 'HNE'
$ws.Cells.Item(1349, 5).Value = '105'
$ws.Cells.Item(1349, 6).Value = 'Return to HNES 003'

# Row 1350
$ws.Cells.Item(1350, 1).Value = 'Pickup PC'
$ws.Cells.Item(1350, 2).Value = 42760
$ws.Cells.Item(1350, 3).Value = '1730'
$ws.Cells.Item(1350, 4).Value = 'HNE'
$ws.Cells.Item(1350, 5).Value = '105'
$ws.Cells.Item(1350, 6).Value = 'Return to HNES 003'

# Row 1351
$ws.Cells.Item(1351, 1).Value = 'Pickup Projector'
$ws.Cells.Item(1351, 2).Value = 42760
$ws.Cells.Item(1351, 3).Value = '1715'
$ws.Cells.Item(1351, 4).Value = 'ATK'
$ws.Cells.Item(1351, 5).Value = '005'
$ws.Cells.Item(1351, 6).Value = 'return projector to DB 0003 storeroom'

# Row 1352
$ws.Cells.Item(1352, 1).Value = 'Pickup PC'
$ws.Cells.Item(1352, 2).Value = 42760
$ws.Cells.Item(1352, 3).Value = '1715'
$ws.Cells.Item(1352, 4).Value = 'ATK'
$ws.Cells.Item(1352, 5).Value = '005'
$ws.Cells.Item(1352, 6).Value = 'return equipment to DB 0003 and plug in for updates'

# Row 1353
$ws.Cells.Item(1353, 1).Value = 'AV Shutdown'
$ws.Cells.Item(1353, 2).Value = 42760
$ws.Cells.Item(1353, 3).Value = '1600'
$ws.Cells.Item(1353, 4).Value = 'KT'
$ws.Cells.Item(1353, 5).Value = '519'

# Row 1354
$ws.Cells.Item(1354, 1).Value = 'Pickup Mic'
$ws.Cells.Item(1354, 2).Value = 42760
$ws.Cells.Item(1354, 3).Value = '1630'
$ws.Cells.Item(1354, 4).Value = 'YL'
$ws.Cells.Item(1354, 5).Value = '242'
$ws.Cells.Item(1354, 6).Value = 'Pick up 4 desk mics, stands and all mic cables and ac cords. Pick up mixer - return all equipment to YL 203C storeroom.'
$ws.Rows.Item(1354).RowHeight = 30

# Row 1355
$ws.Cells.Item(1355, 1).Value = 'Demo'
$ws.Cells.Item(1355, 2).Value = 42760
$ws.Cells.Item(1355, 3).Value = '1630'
$ws.Cells.Item(1355, 4).Value = 'OSG'
$ws.Cells.Item(1355, 5).Value = '2003'

# Row 1356
$ws.Cells.Item(1356, 1).Value = 'SCLD Student Event'
$ws.Cells.Item(1356, 2).Value = 42760
$ws.Cells.Item(1356, 3).Value = '1600'
$ws.Cells.Item(1356, 4).Value = 'FC'
$ws.Cells.Item(1356, 5).Value = '152 - Assembly Hall'
$ws.Cells.Item(1356, 6).Value = 'Student group here, please turn on PC, projector - provide wireless keyboard from FC 156A'
$ws.Rows.Item(1356).RowHeight = 30

# Row 1357
$ws.Cells.Item(1357, 1).Value = 'SCLD Student Logout'
$ws.Cells.Item(1357, 2).Value = 42760
$ws.Cells.Item(1357, 3).Value = '2000'
$ws.Cells.Item(1357, 4).Value = 'FC'
$ws.Cells.Item(1357, 5).Value = '152 - Assembly Hall'
$ws.Cells.Item(1357, 6).Value = 'Turn off projector and PC, return wireless keyboard and projector remote to FC 156A'
$ws.Rows.Item(1357).RowHeight = 30

# Row 1362
$ws.Cells.Item(1362, 1).Value = 'Pickup PC'
$ws.Cells.Item(1362, 2).Value = 42761
$ws.Cells.Item(1362, 3).Value = '1730'
$ws.Cells.Item(1362, 4).Value = 'OSG'
$ws.Cells.Item(1362, 5).Value = '4034'
$ws.Cells.Item(1362, 6).Value = 'Return to OSG 1014L'

# Row 1367
$ws.Cells.Item(1367, 1).Value = 'Pickup Small PA'
$ws.Cells.Item(1367, 2).Value = 42765
$ws.Cells.Item(1367, 3).Value = '1730'
$ws.Cells.Item(1367, 4).Value = 'SSB'
$ws.Cells.Item(1367, 5).Value = 'S126'
$ws.Cells.Item(1367, 6).Value = 'Return Small PA and neck mic to SSB N103'

# Row 1368
$ws.Cells.Item(1368, 1).Value = 'Pickup Skype Kit'
$ws.Cells.Item(1368, 2).Value = 42765
$ws.Cells.Item(1368, 3).Value = '1730'
$ws.Cells.Item(1368, 4).Value = 'OSG'
$ws.Cells.Item(1368, 5).Value = '2009'
$ws.Cells.Item(1368, 6).Value = 'Return Skype kit to OSG 1014L'

# Row 1373
$ws.Cells.Item(1373, 1).Value = 'Pickup Projector'
$ws.Cells.Item(1373, 2).Value = 42766
$ws.Cells.Item(1373, 3).Value = '1730'
$ws.Cells.Item(1373, 4).Value = 'SSB'
$ws.Cells.Item(1373, 5).Value = 'W356'
$ws.Cells.Item(1373, 6).Value = 'Returjn doc camera on cart to SSB N103'

# Row 1374
$ws.Cells.Item(1374, 1).Value = 'Demo'
$ws.Cells.Item(1374, 2).Value = 42766
$ws.Cells.Item(1374, 3).Value = '1600'
$ws.Cells.Item(1374, 4).Value = 'DB'
$ws.Cells.Item(1374, 5).Value = '2008'

# Row 1375
$ws.Cells.Item(1375, 1).Value = 'Demo'
$ws.Cells.Item(1375, 2).Value = 42766
$ws.Cells.Item(1375, 3).Value = '1600'
$ws.Cells.Item(1375, 4).Value = 'HNE'
$ws.Cells.Item(1375, 5).Value = '104'

# Row 1376
$ws.Cells.Item(1376, 1).Value = 'Demo'
$ws.Cells.Item(1376, 2).Value = 42766
$ws.Cells.Item(1376, 3).Value = '1545'
$ws.Cells.Item(1376, 4).Value = 'HNE'
$ws.Cells.Item(1376, 5).Value = '032'
$ws.Cells.Item(1376, 6).Value = 'demo neck mic'

# Row 1377
$ws.Cells.Item(1377, 1).Value = 'Setup Skype Kit'
$ws.Cells.Item(1377, 2).Value = 42766
$ws.Cells.Item(1377, 3).Value = '1615'
$ws.Cells.Item(1377, 4).Value = 'OSG'
$ws.Cells.Item(1377, 5).Value = '1001'
$ws.Cells.Item(1377, 6).Value = 'Skype kit from OSG 1014L'

# Row 1378
$ws.Cells.Item(1378, 1).Value = 'Pickup Skype Kit'
$ws.Cells.Item(1378, 2).Value = 42766
$ws.Cells.Item(1378, 3).Value = '1830'
$ws.Cells.Item(1378, 4).Value = 'OSG'
$ws.Cells.Item(1378, 5).Value = '1001'
$ws.Cells.Item(1378, 6).Value = 'Return Skype kit to OSG 1014L'

# Remove the list-data-validation from D1333 (blank gap row) so the
# validation range splits exactly like the source: D1074:D1332 D1334:D1048576
$ws.Range("D1333").Validation.Delete()

# Leave the sheet selection/view matching the end of the edit session.
$ws.Activate()
$ws.Range("F1386").Select()
